$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 316.03
$ws.Range("I15").Value = 316.03
$ws.Range("K15").Value = 948.0899999999999
$ws.Range("M15").Value = -779.0899999999999
# Row 17
$ws.Range("H17").Value = 1802935.1
$ws.Range("J17").Value = 1802935.1
$ws.Range("L17").Value = 5408805.300000001
$ws.Range("N17").Value = -5409141.300000001
# Row 80
$ws.Range("H80").Value = 465.92
$ws.Range("I80").Value = 373.2
$ws.Range("J80").Value = 605
$ws.Range("K80").Value = 1119.6
$ws.Range("L80").Value = 1815
$ws.Range("M80").Value = -121.5999999999999
$ws.Range("N80").Value = -3811
# Row 83
$ws.Range("H83").Value = 465.92
$ws.Range("I83").Value = 373.2
$ws.Range("J83").Value = 605
$ws.Range("K83").Value = 3358.8
$ws.Range("L83").Value = 5445
$ws.Range("M83").Value = 1633.2
$ws.Range("N83").Value = -15429
# Row 111
$ws.Range("H111").Value = 68682
$ws.Range("I111").Value = 2172.9
$ws.Range("K111").Value = 6518.700000000001
$ws.Range("M111").Value = -3451.700000000001
# Row 129
$ws.Range("H129").Value = 962.5323
$ws.Range("I129").Value = 769.25
$ws.Range("K129").Value = 2307.75
$ws.Range("M129").Value = 2692.25
# Row 132
$ws.Range("H132").Value = 1533.9487
$ws.Range("I132").Value = 1129.8518
$ws.Range("J132").Value = 2443.1667
$ws.Range("K132").Value = 3389.5554
$ws.Range("L132").Value = 7329.500100000001
$ws.Range("M132").Value = -859.5553999999997
$ws.Range("N132").Value = -12389.5001
# Row 135
$ws.Range("H135").Value = 1797.375
$ws.Range("I135").Value = 1122.4333
$ws.Range("J135").Value = 3822.2
$ws.Range("K135").Value = 10101.8997
$ws.Range("L135").Value = 34399.8
$ws.Range("M135").Value = -7566.8997
$ws.Range("N135").Value = -39469.8
# Row 137
$ws.Range("H137").Value = 1218.3939
$ws.Range("I137").Value = 986.9820999999999
$ws.Range("J137").Value = 2514.3
$ws.Range("K137").Value = 2960.9463
$ws.Range("L137").Value = 7542.900000000001
$ws.Range("M137").Value = -410.9462999999996
$ws.Range("N137").Value = -12642.9
# Row 141
$ws.Range("H141").Value = 1766.1212
$ws.Range("I141").Value = 1428.25
$ws.Range("K141").Value = 4284.75
$ws.Range("M141").Value = 895.25

$ws = $wb.Worksheets.Item("ARM")
# Row 26
$ws.Range("H26").Value = 4314.25
$ws.Range("I26").Value = 3085.6667
$ws.Range("K26").Value = 3085.6667
$ws.Range("M26").Value = -2755.6667
# Row 32
$ws.Range("H32").Value = 884.89
$ws.Range("I32").Value = 757.0341
$ws.Range("J32").Value = 1822.5
$ws.Range("K32").Value = 757.0341
$ws.Range("L32").Value = 1822.5
$ws.Range("M32").Value = -470.0341
$ws.Range("N32").Value = -2396.5
# Row 63
$ws.Range("H63").Value = 200004460
$ws.Range("I63").Value = 250005280
$ws.Range("J63").Value = 1250
$ws.Range("K63").Value = 250005280
$ws.Range("L63").Value = 1250
$ws.Range("M63").Value = -250004594
$ws.Range("N63").Value = -2622
# Row 66
$ws.Range("H66").Value = 200004460
$ws.Range("I66").Value = 250005280
$ws.Range("J66").Value = 1250
$ws.Range("K66").Value = 1250026400
$ws.Range("L66").Value = 6250
$ws.Range("M66").Value = -1250022968
$ws.Range("N66").Value = -13114
# Row 102
$ws.Range("H102").Value = 3088354.8
$ws.Range("I102").Value = 3088354.8
$ws.Range("K102").Value = 3088354.8
$ws.Range("M102").Value = -3086732.8
# Row 122
$ws.Range("H122").Value = 2334513.5
$ws.Range("I122").Value = 3667821.2
$ws.Range("J122").Value = 1225
$ws.Range("K122").Value = 11003463.6
$ws.Range("L122").Value = 3675
$ws.Range("M122").Value = -11001013.6
$ws.Range("N122").Value = -8575

$ws = $wb.Worksheets.Item("BSM")
# Row 61
$ws.Range("H61").Value = 25000
$ws.Range("J61").Value = 25000
$ws.Range("L61").Value = 25000
$ws.Range("N61").Value = -25626
# Row 94
$ws.Range("H94").Value = 1602.9375
$ws.Range("I94").Value = 547.8333
$ws.Range("J94").Value = 2236
$ws.Range("K94").Value = 547.8333
$ws.Range("L94").Value = 2236
$ws.Range("M94").Value = -96.83330000000001
$ws.Range("N94").Value = -3138
# Row 134
$ws.Range("H134").Value = 3635.6938
$ws.Range("I134").Value = 4024.25
$ws.Range("J134").Value = 2559.6924
$ws.Range("K134").Value = 12072.75
$ws.Range("L134").Value = 7679.0772
$ws.Range("M134").Value = -9537.75
$ws.Range("N134").Value = -12749.0772

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 6246.224
$ws.Range("I31").Value = 1411.8649
$ws.Range("J31").Value = 12208.6
$ws.Range("K31").Value = 1411.8649
$ws.Range("L31").Value = 12208.6
$ws.Range("M31").Value = -1116.8649
$ws.Range("N31").Value = -12798.6
# Row 34
$ws.Range("H34").Value = 6246.224
$ws.Range("I34").Value = 1411.8649
$ws.Range("J34").Value = 12208.6
$ws.Range("K34").Value = 1411.8649
$ws.Range("L34").Value = 12208.6
$ws.Range("M34").Value = -1209.8649
$ws.Range("N34").Value = -12612.6
# Row 122
$ws.Range("H122").Value = 2532
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 2684
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 8052
$ws.Range("M122").Value = -3550
$ws.Range("N122").Value = -12952
# Row 132
$ws.Range("H132").Value = 1648.4237
$ws.Range("I132").Value = 1342.2195
$ws.Range("J132").Value = 2345.889
$ws.Range("K132").Value = 4026.6585
$ws.Range("L132").Value = 7037.667
$ws.Range("M132").Value = -1496.6585
$ws.Range("N132").Value = -12097.667
# Row 141
$ws.Range("H141").Value = 34880.23
$ws.Range("J141").Value = 34880.23
$ws.Range("L141").Value = 34880.23
$ws.Range("N141").Value = -45240.23

$ws = $wb.Worksheets.Item("CUL")
# Row 87
$ws.Range("H87").Value = 4916.6665
$ws.Range("I87").Value = 4975
$ws.Range("J87").Value = 4800
$ws.Range("K87").Value = 14925
$ws.Range("L87").Value = 14400
$ws.Range("M87").Value = -13677
$ws.Range("N87").Value = -16896
# Row 90
$ws.Range("H90").Value = 4916.6665
$ws.Range("I90").Value = 4975
$ws.Range("J90").Value = 4800
$ws.Range("K90").Value = 44775
$ws.Range("L90").Value = 43200
$ws.Range("M90").Value = -38535
$ws.Range("N90").Value = -55680

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 2402443.8
$ws.Range("I122").Value = 4052735.2
$ws.Range("J122").Value = 2019.3636
$ws.Range("K122").Value = 12158205.6
$ws.Range("L122").Value = 6058.0908
$ws.Range("M122").Value = -12155755.6
$ws.Range("N122").Value = -10958.0908
# Row 132
$ws.Range("H132").Value = 1939.3556
$ws.Range("I132").Value = 1879.8334
$ws.Range("K132").Value = 5639.5002
$ws.Range("M132").Value = -3109.5002

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 700
$ws.Range("I16").Value = 400
$ws.Range("J16").Value = 850
$ws.Range("K16").Value = 400
$ws.Range("L16").Value = 850
$ws.Range("M16").Value = -230
$ws.Range("N16").Value = -1190
# Row 22
$ws.Range("H22").Value = 5436303.5
$ws.Range("I22").Value = 25000260
$ws.Range("J22").Value = 1871.2222
$ws.Range("K22").Value = 25000260
$ws.Range("L22").Value = 1871.2222
$ws.Range("M22").Value = -24999965
$ws.Range("N22").Value = -2461.2222
# Row 27
$ws.Range("H27").Value = 5436303.5
$ws.Range("I27").Value = 25000260
$ws.Range("J27").Value = 1871.2222
$ws.Range("K27").Value = 25000260
$ws.Range("L27").Value = 1871.2222
$ws.Range("M27").Value = -25000153
$ws.Range("N27").Value = -2085.2222
# Row 40
$ws.Range("H40").Value = 27030148
$ws.Range("I40").Value = 34485030
$ws.Range("J40").Value = 6195
$ws.Range("K40").Value = 34485030
$ws.Range("L40").Value = 6195
$ws.Range("M40").Value = -34484894
$ws.Range("N40").Value = -6467
# Row 82
$ws.Range("H82").Value = 74810.42999999999
$ws.Range("I82").Value = 1419.4
$ws.Range("J82").Value = 115583.22
$ws.Range("K82").Value = 1419.4
$ws.Range("L82").Value = 115583.22
$ws.Range("M82").Value = -1058.4
$ws.Range("N82").Value = -116305.22
# Row 85
$ws.Range("H85").Value = 74810.42999999999
$ws.Range("I85").Value = 1419.4
$ws.Range("J85").Value = 115583.22
$ws.Range("K85").Value = 1419.4
$ws.Range("L85").Value = 115583.22
$ws.Range("M85").Value = -171.4000000000001
$ws.Range("N85").Value = -118079.22
# Row 122
$ws.Range("H122").Value = 2555792
$ws.Range("I122").Value = 3406656
$ws.Range("J122").Value = 3200
$ws.Range("K122").Value = 10219968
$ws.Range("L122").Value = 9600
$ws.Range("M122").Value = -10217518
$ws.Range("N122").Value = -14500

$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 600.75
$ws.Range("I100").Value = 401.5
$ws.Range("K100").Value = 803
$ws.Range("M100").Value = -262
# Row 122
$ws.Range("H122").Value = 2606.158
$ws.Range("I122").Value = 2527.8
$ws.Range("K122").Value = 7583.400000000001
$ws.Range("M122").Value = -5133.400000000001
# Row 132
$ws.Range("H132").Value = 39856.96
$ws.Range("I132").Value = 63687.062
$ws.Range("J132").Value = 1728.8
$ws.Range("K132").Value = 191061.186
$ws.Range("L132").Value = 5186.4
$ws.Range("M132").Value = -188531.186
$ws.Range("N132").Value = -10246.4
# Row 136
$ws.Range("H136").Value = 6331019.5
$ws.Range("I136").Value = 1993.3889
$ws.Range("K136").Value = 5980.1667
$ws.Range("M136").Value = -3430.1667
